$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.04187148740654
$ws.Range("D2").Value = 1.051201815403074
$ws.Range("E2").Value = 1.039821650265975
$ws.Range("F2").Value = 1.058114308825262
$ws.Range("I2").Value = 1.02359499962809
$ws.Range("J2").Value = 1.046950360256459
$ws.Range("K2").Value = 1.053953972737855
$ws.Range("L2").Value = 1.042605748426806
$ws.Range("M2").Value = 1.06084742817026
$ws.Range("N2").Value = 1.018867705991366

$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.044480530114546
$ws.Range("D3").Value = 1.053794344556942
$ws.Range("E3").Value = 1.042104672355415
$ws.Range("F3").Value = 1.060870762483992
$ws.Range("I3").Value = 1.023504579208684
$ws.Range("J3").Value = 1.049197364533617
$ws.Range("K3").Value = 1.056354528623648
$ws.Range("L3").Value = 1.044695193457337
$ws.Range("M3").Value = 1.063412932042585
$ws.Range("N3").Value = 1.019675759831763

$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.046159367466952
$ws.Range("D4").Value = 1.055462855549148
$ws.Range("E4").Value = 1.043573327982253
$ws.Range("F4").Value = 1.062645297380346
$ws.Range("I4").Value = 1.023443335729026
$ws.Range("J4").Value = 1.050642091253776
$ws.Range("K4").Value = 1.05789857700645
$ws.Range("L4").Value = 1.046038290180849
$ws.Range("M4").Value = 1.065063695222776
$ws.Range("N4").Value = 1.020193834311687

$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.046862965162818
$ws.Range("D5").Value = 1.056162196342181
$ws.Range("E5").Value = 1.044188741922257
$ws.Range("F5").Value = 1.063389202993886
$ws.Range("I5").Value = 1.023416932628352
$ws.Range("J5").Value = 1.051247297617153
$ws.Range("K5").Value = 1.058545532444821
$ws.Range("L5").Value = 1.046600844292042
$ws.Range("M5").Value = 1.065755512594656
$ws.Range("N5").Value = 1.020410505283552

$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.046980975933702
$ws.Range("D6").Value = 1.056279497256249
$ws.Range("E6").Value = 1.044291956351954
$ws.Range("F6").Value = 1.063513986188949
$ws.Range("I6").Value = 1.023412460900047
$ws.Range("J6").Value = 1.051348789515656
$ws.Range("K6").Value = 1.058654033945397
$ws.Range("L6").Value = 1.046695178856999
$ws.Range("M6").Value = 1.065871546647893
$ws.Range("N6").Value = 1.020446819728551

$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.0461687774792
$ws.Range("D7").Value = 1.055472208351387
$ws.Range("E7").Value = 1.043581558996655
$ws.Range("F7").Value = 1.062655245685326
$ws.Range("I7").Value = 1.02344298551034
$ws.Range("J7").Value = 1.050650186453369
$ws.Range("K7").Value = 1.057907230077869
$ws.Range("L7").Value = 1.046045815177342
$ws.Range("M7").Value = 1.065072947741051
$ws.Range("N7").Value = 1.020196733884585

$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.042755207016752
$ws.Range("D8").Value = 1.052079878841492
$ws.Range("E8").Value = 1.040595023663328
$ws.Range("F8").Value = 1.059047782971711
$ws.Range("I8").Value = 1.023565007130014
$ws.Range("J8").Value = 1.047711691033278
$ws.Range("K8").Value = 1.054767206369114
$ws.Range("L8").Value = 1.043313762384146
$ws.Range("M8").Value = 1.061716412354042
$ws.Range("N8").Value = 1.019141795176683

$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.036665337769893
$ws.Range("D9").Value = 1.046030263772806
$ws.Range("E9").Value = 1.035263993507622
$ws.Range("F9").Value = 1.052618565370167
$ws.Range("I9").Value = 1.023759139049229
$ws.Range("J9").Value = 1.042460483694758
$ws.Range("K9").Value = 1.049160464999226
$ws.Range("L9").Value = 1.038428974103899
$ws.Range("M9").Value = 1.055727867516465
$ws.Range("N9").Value = 1.017245289056235

$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.032551179727448
$ws.Range("D10").Value = 1.041944948051617
$ws.Range("E10").Value = 1.031660577821856
$ws.Range("F10").Value = 1.048279616504079
$ws.Range("I10").Value = 1.023874603872166
$ws.Range("J10").Value = 1.038906952928114
$ws.Range("K10").Value = 1.045369439635253
$ws.Range("L10").Value = 1.035121752625018
$ws.Range("M10").Value = 1.051681907630736
$ws.Range("N10").Value = 1.015954440070691

$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.030755928452849
$ws.Range("D11").Value = 1.040162685939144
$ws.Range("E11").Value = 1.03008776563715
$ws.Range("F11").Value = 1.046387344262862
$ws.Range("I11").Value = 1.023921309042974
$ws.Range("J11").Value = 1.037354932637964
$ws.Range("K11").Value = 1.043714426013898
$ws.Range("L11").Value = 1.033676926104313
$ws.Range("M11").Value = 1.049916362575611
$ws.Range("N11").Value = 1.015388915221136

$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.030086940443625
$ws.Range("D12").Value = 1.039498599960276
$ws.Range("E12").Value = 1.029501605844354
$ws.Range("F12").Value = 1.045682363410549
$ws.Range("I12").Value = 1.023938164268118
$ws.Range("J12").Value = 1.036776372736509
$ws.Range("K12").Value = 1.043097582477152
$ws.Range("L12").Value = 1.033138268930366
$ws.Range("M12").Value = 1.049258437118932
$ws.Range("N12").Value = 1.015177840786768

$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.030230539311419
$ws.Range("D13").Value = 1.039641143794275
$ws.Range("E13").Value = 1.02962742833976
$ws.Range("F13").Value = 1.045833680862914
$ws.Range("I13").Value = 1.023934571058996
$ws.Range("J13").Value = 1.03690057062913
$ws.Range("K13").Value = 1.043229993638122
$ws.Range("L13").Value = 1.03325390357013
$ws.Range("M13").Value = 1.049399661718104
$ws.Range("N13").Value = 1.01522316323238

$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.03070067405705
$ws.Range("D14").Value = 1.040107835157795
$ws.Range("E14").Value = 1.030039353609649
$ws.Range("F14").Value = 1.046329113781361
$ws.Range("I14").Value = 1.023922712353277
$ws.Range("J14").Value = 1.037307151341734
$ws.Range("K14").Value = 1.043663480760194
$ws.Range("L14").Value = 1.033632441392224
$ws.Range("M14").Value = 1.049862022024542
$ws.Range("N14").Value = 1.015371488560735

$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.030990052053565
$ws.Range("D15").Value = 1.040395101877316
$ws.Range("E15").Value = 1.030292894265003
$ws.Range("F15").Value = 1.04663408472021
$ws.Range("I15").Value = 1.023915340507297
$ws.Range("J15").Value = 1.037557382670714
$ws.Range("K15").Value = 1.043930286279884
$ws.Range("L15").Value = 1.033865406144137
$ws.Range("M15").Value = 1.050146613745547
$ws.Range("N15").Value = 1.015462741658496

$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.032670027291302
$ws.Range("D16").Value = 1.042062944286947
$ws.Range("E16").Value = 1.0317646908681
$ws.Range("F16").Value = 1.048404909474329
$ws.Range("I16").Value = 1.023871434968932
$ws.Range("J16").Value = 1.039009668907735
$ws.Range("K16").Value = 1.045478987321475
$ws.Range("L16").Value = 1.035217366315297
$ws.Range("M16").Value = 1.053432480014797
$ws.Range("N16").Value = 1.015991831371257

$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.033720082961782
$ws.Range("D17").Value = 1.043105524765031
$ws.Range("E17").Value = 1.0326845134386
$ws.Range("F17").Value = 1.049512035391184
$ws.Range("I17").Value = 1.023843013862359
$ws.Range("J17").Value = 1.039917035496744
$ws.Range("K17").Value = 1.046446787790026
$ws.Range("L17").Value = 1.036061949156615
$ws.Range("M17").Value = 1.052831452570965
$ws.Range("N17").Value = 1.016321936643824

$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.034331237235226
$ws.Range("D18").Value = 1.043712367005878
$ws.Range("E18").Value = 1.033219828073061
$ws.Range("F18").Value = 1.050156508315902
$ws.Range("I18").Value = 1.023826118601775
$ws.Range("J18").Value = 1.040445006863655
$ws.Range("K18").Value = 1.047009994367402
$ws.Range("L18").Value = 1.036553351603943
$ws.Range("M18").Value = 1.053432480014797
$ws.Range("N18").Value = 1.016513848429056

$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.034539402610253
$ws.Range("D19").Value = 1.043919070187845
$ws.Range("E19").Value = 1.033402154631214
$ws.Range("F19").Value = 1.050376039737976
$ws.Range("I19").Value = 1.023820303822167
$ws.Range("J19").Value = 1.040624816393087
$ws.Range("K19").Value = 1.047201815788179
$ws.Range("L19").Value = 1.036720700708765
$ws.Range("M19").Value = 1.053637195203281
$ws.Range("N19").Value = 1.016579178761367

$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.033607559590002
$ws.Range("D20").Value = 1.042993798418151
$ws.Range("E20").Value = 1.032585950028397
$ws.Range("F20").Value = 1.049393385779098
$ws.Range("I20").Value = 1.023846096030143
$ws.Range("J20").Value = 1.039819816565332
$ws.Range("K20").Value = 1.046343086438366
$ws.Range("L20").Value = 1.035971460935525
$ws.Range("M20").Value = 1.052720793296082
$ws.Range("N20").Value = 1.016286585128486

$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.0305622911383
$ws.Range("D21").Value = 1.039970464078725
$ws.Range("E21").Value = 1.029918106193659
$ws.Range("F21").Value = 1.046183279955623
$ws.Range("I21").Value = 1.02392621804881
$ws.Range("J21").Value = 1.037187481177586
$ws.Range("K21").Value = 1.043535888114371
$ws.Range("L21").Value = 1.033521026740185
$ws.Range("M21").Value = 1.049725927686525
$ws.Range("N21").Value = 1.015327838624565

$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.028635119882757
$ws.Range("D22").Value = 1.038057531862724
$ws.Range("E22").Value = 1.028229425464243
$ws.Range("F22").Value = 1.04415272847382
$ws.Range("I22").Value = 1.023973741382756
$ws.Range("J22").Value = 1.035520412689214
$ws.Range("K22").Value = 1.041758716308536
$ws.Range("L22").Value = 1.031968827991333
$ws.Range("M22").Value = 1.047830611641984
$ws.Range("N22").Value = 1.0147191610485

$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.029657960769316
$ws.Range("D23").Value = 1.039072781158898
$ws.Range("E23").Value = 1.029125721426758
$ws.Range("F23").Value = 1.045230349335695
$ws.Range("I23").Value = 1.023948818288664
$ws.Range("J23").Value = 1.036405319775497
$ws.Range("K23").Value = 1.042702007709648
$ws.Range("L23").Value = 1.032792791232974
$ws.Range("M23").Value = 1.04883654919575
$ws.Range("N23").Value = 1.015042397980009

$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.033658408149572
$ws.Range("D24").Value = 1.043044286687922
$ws.Range("E24").Value = 1.032630490302435
$ws.Range("F24").Value = 1.049447002428487
$ws.Range("I24").Value = 1.023844704313072
$ws.Range("J24").Value = 1.039863749567876
$ws.Range("K24").Value = 1.04638994861327
$ws.Range("L24").Value = 1.036012352460081
$ws.Range("M24").Value = 1.052770799498561
$ws.Range("N24").Value = 1.016302560911394

$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.03824897671035
$ws.Range("D25").Value = 1.047603150431547
$ws.Range("E25").Value = 1.036650642989912
$ws.Range("F25").Value = 1.054289671377581
$ws.Range("I25").Value = 1.02371141742819
$ws.Range("J25").Value = 1.043827077826951
$ws.Range("K25").Value = 1.050619041999529
$ws.Range("L25").Value = 1.039700501341906
$ws.Range("M25").Value = 1.057285204556185
$ws.Range("N25").Value = 1.017740160577004
